$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Result")

# Row 2 (A2:O2)
$ws.Range("A2").Value = 0.0790000000000077
$ws.Range("B2").Value = 0.334000000000003
$ws.Range("C2").Value = 0.005
$ws.Range("D2").Value = 0.022
$ws.Range("E2").Value = 0.02
$ws.Range("F2").Value = 0.0
$ws.Range("G2").Value = -0.007
$ws.Range("H2").Value = 0.006
$ws.Range("I2").Value = 0.004
$ws.Range("J2").Value = 150.116
$ws.Range("K2").Value = 170.208
$ws.Range("L2").Value = 251.204
$ws.Range("M2").Value = 222.473
$ws.Range("N2").Value = 148.757
$ws.Range("O2").Value = 130.171

# Row 3 (A3:O3)
$ws.Range("A3").Value = -0.0679999999999978
$ws.Range("B3").Value = -0.0159999999999911
$ws.Range("C3").Value = -0.013
$ws.Range("D3").Value = -0.004
$ws.Range("E3").Value = 0.005
$ws.Range("F3").Value = -0.015
$ws.Range("G3").Value = -0.022
$ws.Range("H3").Value = -0.009
$ws.Range("I3").Value = -0.011
$ws.Range("J3").Value = 149.653
$ws.Range("K3").Value = 169.63
$ws.Range("L3").Value = 249.712
$ws.Range("M3").Value = 222.105
$ws.Range("N3").Value = 148.437
$ws.Range("O3").Value = 129.832
